$d = $word.ActiveDocument

$pairs = @(
    @("2024-11-28 Thursday", "2024-11-29 Friday"),
    @("598÷3=", "947÷9="),
    @("612÷4=", "432÷9="),
    @("306÷3=", "170÷2="),
    @("770÷8=", "418÷5="),
    @("761÷4=", "926÷9="),
    @("316÷9=", "600÷3="),
    @("118÷7=", "388÷7="),
    @("536÷7=", "920÷7="),
    @("886÷2=", "462÷8="),
    @("452÷3=", "118÷4="),
    @("920÷3=", "471÷6="),
    @("179÷9=", "950÷8="),
    @("215÷8=", "366÷9="),
    @("162÷6=", "649÷9="),
    @("764÷8=", "444÷8="),
    @("247÷4=", "751÷4="),
    @("317÷3=", "273÷9="),
    @("579÷5=", "204÷7="),
    @("516÷6=", "211÷2="),
    @("724÷9=", "570÷7="),
    @("496÷7=", "928÷3="),
    @("621÷5=", "825÷7="),
    @("914÷7=", "329÷5="),
    @("823÷2=", "394÷2="),
    @("210÷3=", "793÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
